$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value2 = 44176
$ws.Range("M2").Value2 = 250

# Row 3
$ws.Range("D3").Value2 = 44208
$ws.Range("M3").Value2 = 210
$ws.Range("Q3").Value = "$/caja 14 kilos empedrada"

# Row 4
$ws.Range("D4").Value2 = 44351
$ws.Range("M4").Value2 = 300
$ws.Range("N4").Value2 = 10000
$ws.Range("O4").Value2 = 10000
$ws.Range("P4").Value2 = 10000
$ws.Range("S4").Value2 = 714

# Row 5
$ws.Range("D5").Value2 = 44400
$ws.Range("M5").Value2 = 100
$ws.Range("N5").Value2 = 10000
$ws.Range("O5").Value2 = 10000
$ws.Range("P5").Value2 = 10000
$ws.Range("S5").Value2 = 714

# Row 6
$ws.Range("D6").Value2 = 44309
$ws.Range("M6").Value2 = 300
$ws.Range("N6").Value2 = 7000
$ws.Range("O6").Value2 = 7000
$ws.Range("P6").Value2 = 7000
$ws.Range("S6").Value2 = 500

# Row 8
$ws.Range("D8").Value2 = 44491
$ws.Range("M8").Value2 = 180
$ws.Range("N8").Value2 = 9000
$ws.Range("O8").Value2 = 9000
$ws.Range("P8").Value2 = 9000
$ws.Range("S8").Value2 = 643

# Row 9
$ws.Range("D9").Value2 = 44397
$ws.Range("M9").Value2 = 60
$ws.Range("N9").Value2 = 11000
$ws.Range("O9").Value2 = 11000
$ws.Range("P9").Value2 = 11000
$ws.Range("Q9").Value = "$/caja 14 kilos"
$ws.Range("S9").Value2 = 786
